# "Fix specs for DateTime parsing"
#
# The fixture's data table is trimmed down to its essential 2 columns x 5
# rows (header + blank + Date/Datetime/Datetime rows), and the date/datetime
# sample cells get corrected number formats so the DateTime-parsing specs
# exercise the right display strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to carry 5 extra "filler" columns (C:G) and 18 extra empty
# rows (6:23) that mirrored the formatting of the real data region but held
# no content. Drop them so the sheet shrinks to the real A1:B5 table.
$ws.Range("C:G").Delete()
$ws.Range("6:23").Delete()

# Correct the number formats backing the Date / Datetime sample cells.
$ws.Range("B3").NumberFormat = "mm-dd-yy"
$ws.Range("B4").NumberFormat = "[$-409]d-m-yy h:mm AM/PM;@"
$ws.Range("B5").NumberFormat = "d-m-yy h:mm;@"
